$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("datos")

# Rename province labels in column A to the Ley 39/2015 / INE "co-official
# name first" ordering used in the updated dataset.
$ws.Range("A3").Value  = "Alacant/Alicante"
$ws.Range("A9").Value  = "Illes Balears"
$ws.Range("A16").Value = "Castelló/Castellón"
$ws.Range("A19").Value = "A Coruña"
$ws.Range("A37").Value = "Las Palmas"
$ws.Range("A39").Value = "La Rioja"
$ws.Range("A48").Value = "València/Valencia"

# Add a bit of formatting on J5 (Courier New, dark grey) left over from the
# author poking around the sheet.
$ws.Range("J5").Font.Name = "Courier New"
$ws.Range("J5").Font.Color = 2171169

$ws.Range("J5").Select()
